# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / Leve profit columns (H-N) on several
# sheets to reflect a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1357.7368
$ws.Range("I28").Value = 2160
$ws.Range("J28").Value = 1071.2142
$ws.Range("K28").Value = 2160
$ws.Range("L28").Value = 1071.2142
$ws.Range("M28").Value = -1675
$ws.Range("N28").Value = -2041.2142

$ws.Range("H62").Value = 2417.5356
$ws.Range("I62").Value = 2309.611
$ws.Range("J62").Value = 2611.8
$ws.Range("K62").Value = 2309.611
$ws.Range("L62").Value = 2611.8
$ws.Range("M62").Value = -1685.611
$ws.Range("N62").Value = -3859.8

$ws.Range("H65").Value = 2417.5356
$ws.Range("I65").Value = 2309.611
$ws.Range("J65").Value = 2611.8
$ws.Range("K65").Value = 11548.055
$ws.Range("L65").Value = 13059
$ws.Range("M65").Value = -8428.055
$ws.Range("N65").Value = -19299

$ws.Range("H70").Value = 2422.6667
$ws.Range("I70").Value = 860
$ws.Range("J70").Value = 2564.7273
$ws.Range("K70").Value = 2580
$ws.Range("L70").Value = 7694.1819
$ws.Range("M70").Value = -2310
$ws.Range("N70").Value = -8234.1819

$ws.Range("H73").Value = 2422.6667
$ws.Range("I73").Value = 860
$ws.Range("J73").Value = 2564.7273
$ws.Range("K73").Value = 2580
$ws.Range("L73").Value = 7694.1819
$ws.Range("M73").Value = -1644
$ws.Range("N73").Value = -9566.1819

$ws.Range("H106").Value = 3950.65
$ws.Range("I106").Value = 3632.2632
$ws.Range("K106").Value = 3632.2632
$ws.Range("M106").Value = -3001.2632

$ws.Range("H129").Value = 843.14667
$ws.Range("J129").Value = 969.84485
$ws.Range("L129").Value = 2909.53455
$ws.Range("N129").Value = -12909.53455

$ws.Range("H135").Value = 45091.434
$ws.Range("I135").Value = 59965
$ws.Range("J135").Value = 2949.6667
$ws.Range("K135").Value = 539685
$ws.Range("L135").Value = 26547.0003
$ws.Range("M135").Value = -537150
$ws.Range("N135").Value = -31617.0003

$ws.Range("H138").Value = 995460.3
$ws.Range("I138").Value = 2179
$ws.Range("J138").Value = 1424987.4
$ws.Range("K138").Value = 6537
$ws.Range("L138").Value = 4274962.199999999
$ws.Range("M138").Value = -1397
$ws.Range("N138").Value = -4285242.199999999

$ws.Range("H141").Value = 1625.625
$ws.Range("I141").Value = 1183.3334
$ws.Range("J141").Value = 2952.5
$ws.Range("K141").Value = 3550.0002
$ws.Range("L141").Value = 8857.5
$ws.Range("M141").Value = 1629.9998
$ws.Range("N141").Value = -19217.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16270.271
$ws.Range("I32").Value = 16548.174
$ws.Range("J32").Value = 14672.333
$ws.Range("K32").Value = 16548.174
$ws.Range("L32").Value = 14672.333
$ws.Range("M32").Value = -16261.174
$ws.Range("N32").Value = -15246.333

$ws.Range("H74").Value = 10500951
$ws.Range("I74").Value = 11410060
$ws.Range("K74").Value = 11410060
$ws.Range("M74").Value = -11409186

$ws.Range("H77").Value = 10500951
$ws.Range("I77").Value = 11410060
$ws.Range("K77").Value = 57050300
$ws.Range("M77").Value = -57045932

$ws.Range("H106").Value = 49952.5
$ws.Range("J106").Value = 49952.5
$ws.Range("L106").Value = 49952.5
$ws.Range("N106").Value = -52476.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H76").Value = 39209.332
$ws.Range("J76").Value = 39209.332
$ws.Range("L76").Value = 39209.332
$ws.Range("N76").Value = -39839.332

$ws.Range("H79").Value = 39209.332
$ws.Range("J79").Value = 39209.332
$ws.Range("L79").Value = 39209.332
$ws.Range("N79").Value = -41393.332


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2910.8928
$ws.Range("I31").Value = 1217.6086
$ws.Range("J31").Value = 10700
$ws.Range("K31").Value = 1217.6086
$ws.Range("L31").Value = 10700
$ws.Range("M31").Value = -922.6086
$ws.Range("N31").Value = -11290

$ws.Range("H34").Value = 2910.8928
$ws.Range("I34").Value = 1217.6086
$ws.Range("J34").Value = 10700
$ws.Range("K34").Value = 1217.6086
$ws.Range("L34").Value = 10700
$ws.Range("M34").Value = -1015.6086
$ws.Range("N34").Value = -11104

$ws.Range("H98").Value = 57000
$ws.Range("J98").Value = 57000
$ws.Range("L98").Value = 57000
$ws.Range("N98").Value = -61492

$ws.Range("H105").Value = 1703.1666
$ws.Range("I105").Value = 1303.1818
$ws.Range("K105").Value = 1303.1818
$ws.Range("M105").Value = 443.8181999999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 19818.182
$ws.Range("J94").Value = 19818.182
$ws.Range("L94").Value = 19818.182
$ws.Range("N94").Value = -21170.182

$ws.Range("H100").Value = 40538
$ws.Range("J100").Value = 40538
$ws.Range("L100").Value = 40538
$ws.Range("N100").Value = -42702

$ws.Range("H101").Value = 69580
$ws.Range("J101").Value = 69580
$ws.Range("L101").Value = 69580
$ws.Range("N101").Value = -76070

$ws.Range("H113").Value = 1971.6364
$ws.Range("I113").Value = 1811
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 1811
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 359
$ws.Range("N113").Value = -6740


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 40144
$ws.Range("J76").Value = 40144
$ws.Range("L76").Value = 40144
$ws.Range("N76").Value = -40820

$ws.Range("H79").Value = 40144
$ws.Range("J79").Value = 40144
$ws.Range("L79").Value = 40144
$ws.Range("N79").Value = -42484

$ws.Range("H103").Value = 29355
$ws.Range("J103").Value = 29355
$ws.Range("L103").Value = 29355
$ws.Range("N103").Value = -31699


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 202455.8
$ws.Range("I132").Value = 333866.66
$ws.Range("K132").Value = 1001599.98
$ws.Range("M132").Value = -999069.98

$ws.Range("H136").Value = 44466.74
$ws.Range("I136").Value = 24761.072
$ws.Range("J136").Value = 251376.25
$ws.Range("K136").Value = 74283.216
$ws.Range("L136").Value = 754128.75
$ws.Range("M136").Value = -71733.216
$ws.Range("N136").Value = -759228.75

